$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row updates
$ws.Range("D1").Value = "description"
$ws.Range("H1").Value = "bis"

# Row 2 (item1): location dunetown -> skirmish, add H2 = false
$ws.Range("E2").Value = "skirmish"
$ws.Range("H2").Formula = "'false"

# Row 3 (item2): fill in missing race (G3) and add H3 = false
$ws.Range("G3").Value = "human"
$ws.Range("H3").Formula = "'false"

# Row 4 (new item2 variant)
$ws.Range("A4").Value = "item2"
$ws.Range("B4").Value = "Pulson grenade " + [char]8220 + "Doom D3" + [char]8221
$ws.Range("C4").Value = "purple"
$ws.Range("D4").Value = "trooper,lord commander"
$ws.Range("E4").Value = "skirmish"
$ws.Range("F4").Value = "xenotronics"
$ws.Range("G4").Value = "human"
$ws.Range("H4").Formula = "'false"

# Match the post-edit selection left by the author
$null = $ws.Range("C5").Select()
